# Added filtering options for the Component Analysis
# For each data row, the columns beyond the applicable forecast-filter
# horizon are cleared so only the columns relevant to the configured
# filtering window remain populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @("I2:K2", "G3:K3", "I4:K4", "G5:K5", "I6:K6", "G7:K7", "I8:K8", "G9:K9", "I10:K10", "G11:K11", "I12:K12", "K13", "I14:K14", "K15", "I16:K16", "K17", "J18:K18", "I19:K19", "K21", "J22:K22", "I23:K23", "K25", "J26:K26", "I27:K27", "K29", "J30:K30", "I31:K31", "K33", "J34:K34", "I35:K35", "K37", "J38:K38", "I39:J39")

foreach ($r in $ranges) {
    $ws.Range($r).ClearContents()
}
